# Rename AHB-Diff header columns from generic _old/_new suffixes to the
# concrete format-version suffixes (_FV2310 / _FV2404), freeze the header
# row, and turn the data range into a proper Excel Table (ListObject).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row (row 1, columns A:U) to use the new suffixes.
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Freeze the header row: split after row 1, keep the header visible while
#    scrolling, and select the top-left cell of the scrollable pane.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn A1:U64 into an Excel Table ("Table1") so that the header row gains
#    filter buttons and the sheet carries a tableParts reference.
$range = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
